# qabul.xlsx update — 10:02 time 02.07.2025 date
# 1) Fix a mis-typed applicant name in row 132.
# 2) Append six newly-submitted applicants (rows 146-151).
# 3) Move the on-screen selection down to the newest edited row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Correct the applicant name in A132 -------------------------------
$ws.Range("A132").Value = 'IRSALIYEV TO‘YCHIBOY OTAKUZI O‘G‘LI'

# --- 2) Append the new applicant rows ------------------------------------
# Columns: A F.I.Sh | B Yo'nalish | C Ta'lim tili | D Ta'lim shakli |
#          E Passport | F JSHIR | G Viloyat | H Tuman |
#          I Telegram raqami | J Telefon raqami | K Sana
# Columns E:K hold ID/phone numbers and date-like text that Excel would
# otherwise auto-coerce into numbers/dates, so force them to Text first
# (matches the source data, which stores these as literal strings).

function Add-Applicant {
    param(
        [int]$Row,
        [string]$Fish,
        [string]$Yonalish,
        [string]$TalimTili,
        [string]$TalimShakli,
        [string]$Passport,
        [string]$Jshir,
        [string]$Viloyat,
        [string]$Tuman,
        [string]$Telegram,
        [string]$Telefon,
        [string]$Sana
    )

    $ws.Range("E$($Row):K$($Row)").NumberFormat = "@"

    $ws.Range("A$Row").Value = $Fish
    $ws.Range("B$Row").Value = $Yonalish
    $ws.Range("C$Row").Value = $TalimTili
    $ws.Range("D$Row").Value = $TalimShakli
    $ws.Range("E$Row").Value = $Passport
    $ws.Range("F$Row").Value = $Jshir
    $ws.Range("G$Row").Value = $Viloyat
    $ws.Range("H$Row").Value = $Tuman
    $ws.Range("I$Row").Value = $Telegram
    $ws.Range("J$Row").Value = $Telefon
    $ws.Range("K$Row").Value = $Sana
}

Add-Applicant -Row 146 `
    -Fish "Sa'dinov Husen Islomovich" -Yonalish "Yurisprudensiya" `
    -TalimTili "O'zbek tili" -TalimShakli "Kunduzgi" `
    -Passport "AD1956155" -Jshir "52005056150012" `
    -Viloyat "Samarqand viloyati" -Tuman "Toyloq tumani" `
    -Telegram "998992949844" -Telefon "+998990488733" -Sana "2025-07-01"

Add-Applicant -Row 147 `
    -Fish "Begmatov Hojiakbar Ulugbek ogli" -Yonalish "Yurisprudensiya" `
    -TalimTili "O'zbek tili" -TalimShakli "Kunduzgi" `
    -Passport "AD7482039" -Jshir "51912076580039" `
    -Viloyat "Toshkent shahri" -Tuman "Yashnaobod tumani" `
    -Telegram "998906610323" -Telefon "+998906610323" -Sana "2025-07-01"

Add-Applicant -Row 148 `
    -Fish "Xushboqov Bunyod Toxirovich" -Yonalish "Yurisprudensiya" `
    -TalimTili "O'zbek tili" -TalimShakli "Kunduzgi" `
    -Passport "AD6046888" -Jshir "51202076350025" `
    -Viloyat "Surxondaryo viloyati" -Tuman "Termiz shahri" `
    -Telegram "998994261253" -Telefon "+998997161253" -Sana "2025-07-01"

Add-Applicant -Row 149 `
    -Fish "Mavleeva Elsana Timurovna" -Yonalish "Psixologiya" `
    -TalimTili "Rus tili" -TalimShakli "Kunduzgi" `
    -Passport "AD5687911" -Jshir "60301088660013" `
    -Viloyat "Toshkent shahri" -Tuman "Yashnaobod tumani" `
    -Telegram "998974321132" -Telefon "+998974321132" -Sana "2025-07-01"

Add-Applicant -Row 150 `
    -Fish "Shokirova Dilafruz Eminjon qizi" -Yonalish "Yurisprudensiya" `
    -TalimTili "Rus tili" -TalimShakli "Kunduzgi" `
    -Passport "AD3230069" -Jshir "61003076620056" `
    -Viloyat "Toshkent shahri" -Tuman "Mirzo Ulugʻbek tumani" `
    -Telegram "998900668474" -Telefon "+998331668474" -Sana "2025-07-02"

Add-Applicant -Row 151 `
    -Fish "Tojiddinova Muzayyam" -Yonalish "Yurisprudensiya" `
    -TalimTili "Rus tili" -TalimShakli "Kunduzgi" `
    -Passport "AD6092423" -Jshir "60702085140012" `
    -Viloyat "Toshkent shahri" -Tuman "Yashnaobod tumani" `
    -Telegram "998930639601" -Telefon "+998938390207" -Sana "2025-07-02"

# --- 3) Move the selection to the cell the author last edited ------------
$ws.Application.Goto($ws.Range("A132"), $true)
